# CDF AD & CDF CH Update
# Applies the commit's data corrections to the CDF_AD_sub_1 sheet:
#  - clears screening_encounter_id / last_encounter for rows whose result
#    reverted to "No screening recorded" (and flips numerator/medicaid
#    booleans accordingly)
#  - row 11 (patient 56517) is corrected to a Negative screening result
#  - a couple of stray medicaid flags (H3, H4) are corrected to FALSE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple medicaid-flag corrections ---
$ws.Range("H3").Value = $false
$ws.Range("H4").Value = $false

# --- row 7: becomes "No screening recorded" ---
$ws.Range("D7").Value = ""
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = ""
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = "No screening recorded"
$ws.Range("H7").Value = $false

# --- row 8: becomes "No screening recorded" ---
$ws.Range("D8").Value = ""
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = ""
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = "No screening recorded"
$ws.Range("H8").Value = $false

# --- row 10: becomes "No screening recorded" ---
$ws.Range("D10").Value = ""
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = ""
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = $false
$ws.Range("G10").Value = "No screening recorded"
$ws.Range("H10").Value = $false

# --- row 11: corrected to a Negative screening result ---
# (screening_encounter_id is a text identifier like the other rows, so
# force text formatting before assigning, otherwise Excel would store
# the numeric-looking string as a number)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = ""
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = "Negative screening"
$ws.Range("H11").Value = $true

# --- row 12: becomes "No screening recorded" ---
$ws.Range("D12").Value = ""
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = ""
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = $false
$ws.Range("G12").Value = "No screening recorded"
$ws.Range("H12").Value = $false

# --- row 13: becomes "No screening recorded" ---
$ws.Range("D13").Value = ""
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = ""
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = $false
$ws.Range("G13").Value = "No screening recorded"
$ws.Range("H13").Value = $false

# --- row 14: becomes "No screening recorded" ---
$ws.Range("D14").Value = ""
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = ""
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = $false
$ws.Range("G14").Value = "No screening recorded"
$ws.Range("H14").Value = $true

# --- row 16: becomes "No screening recorded" ---
$ws.Range("D16").Value = ""
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = ""
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = $false
$ws.Range("G16").Value = "No screening recorded"
$ws.Range("H16").Value = $false
